$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:D").Insert()
$ws.Range("D5:D102").Value2 = $ws.Range("E5:E102").Value2
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
Write-Host "done"
